# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Vega Modelo de Temuco" - Platano
# as row 302, shifting the existing rows 302-351 down to 303-352.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 302 (existing rows 302.. move down by one)
$ws.Rows.Item(302).Insert()

# Fill in the new row 302 with the new data record
$ws.Range('A302').Value = 10
$ws.Range('B302').Value = 'Vega Modelo de Temuco'
$ws.Range('C302').Value = 'La Araucanía'
$ws.Range('D302').Value = 44491
$ws.Range('E302').Value = 9
$ws.Range('F302').Value = 'Fruta'
$ws.Range('G302').Value = 100108
$ws.Range('H302').Value = 'Tropicales y subtropicales'
$ws.Range('I302').Value = 100108006
$ws.Range('J302').Value = 'Plátano'
$ws.Range('K302').Value = 'Sin especificar'
$ws.Range('L302').Value = 'Pintón'
$ws.Range('M302').Value = 600
$ws.Range('N302').Value = 25000
$ws.Range('O302').Value = 25000
$ws.Range('P302').Value = 25000
$ws.Range('Q302').Value = '$/caja 20 kilos'
$ws.Range('R302').Value = 'Ecuador'
$ws.Range('S302').Value = 1250
$ws.Range('T302').Value = 20
